$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A:C slightly narrower (41.140625 -> ~40.42578125 char-width units) ---
$ws.Range("A1:C1").EntireColumn.ColumnWidth = 39.67

# --- Row heights: rows 2 and 5-12 become 15pt custom height ---
$ws.Range("A2").EntireRow.RowHeight = 15
$ws.Range("A5:A12").EntireRow.RowHeight = 15

# --- New column M (year 2023 data) ---
# Copy formatting (styles) from column L down into column M first ...
$ws.Range("L3:L12").Copy()
$ws.Range("M3:M12").PasteSpecial(-4122)

# ... then write in the 2023 values (M3, M6 and M9 stay blank, matching L3/L6/L9)
$ws.Range("M4").Value = 2023
$ws.Range("M5").Value = 311.65582791395695
$ws.Range("M7").Value = 119.55977988994496
$ws.Range("M8").Value = 192.09604802401199
$ws.Range("M10").Value = 78.539269634817401
$ws.Range("M11").Value = 60.030015007503756
$ws.Range("M12").Value = 26.013006503251628

# --- Clear the previous explicit cell selection (N5) so the sheet opens at A1 ---
$null = $ws.Range("A1").Select()
